$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New columns "other constellations detected(without degenerate triangles)" (G)
# and "time(seconds) w/ parallelization without degenerate triangles" (J) were
# filled in for rows 8-13. Copy the existing number formatting from the
# neighbouring F/I columns (style index 1) before setting the values.
8..13 | ForEach-Object {
    $r = $_
    $ws.Range("F$r").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
    $ws.Range("I$r").Copy()
    $ws.Range("J$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("G8").Value = 0
$ws.Range("J8").Value = 15

$ws.Range("G9").Value = 0
$ws.Range("J9").Value = 29

$ws.Range("G10").Value = 2
$ws.Range("J10").Value = 50

$ws.Range("G11").Value = 0
$ws.Range("J11").Value = 16

$ws.Range("G12").Value = 1
$ws.Range("J12").Value = 31

$ws.Range("G13").Value = 0
$ws.Range("J13").Value = 15

# Update the saved view state (scroll position + selection) to match.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("J14").Select()
